$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("drop_dups")
$ws.Activate()

# Swap the contents of rows 4 and 5 (A:C) - removing the now-adjacent
# duplicate "I should"/"I should " rows ended up reordered as part of
# the dedup work, so put "austin/ritter" back above "I should/be removed".
$row4 = @($ws.Range("A4").Value(), $ws.Range("B4").Value(), $ws.Range("C4").Value())
$row5 = @($ws.Range("A5").Value(), $ws.Range("B5").Value(), $ws.Range("C5").Value())

$ws.Range("A4").Value = $row5[0]
$ws.Range("B4").Value = $row5[1]
$ws.Range("C4").Value = $row5[2]

$ws.Range("A5").Value = $row4[0]
$ws.Range("B5").Value = $row4[1]
$ws.Range("C5").Value = $row4[2]

$ws.Range("E6").Select()
